$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the invoice-number column (column C) - shifts column D (Fizetendő) into C
$ws.Columns("C").Delete()

$ws.Range("A1:C1").Select()
